# More introduction screens added:
# Duplicate the existing "Rounded Rectangular Callout 3" shape (the
# custom speech-bubble callout already on the slide) to create a third
# callout, then reposition / resize / flip it into its new spot.
#
# Duplicating (rather than re-building the custom geometry from
# scratch) guarantees the new shape keeps byte-identical <a:custGeom>
# path data, fill/line/style and text body to its sibling callout, and
# PowerPoint naturally assigns it the next free shape id (7).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rounded Rectangular Callout 3" is shape index 3 (id=4) on the slide.
$orig = $s.Shapes.Item(3)

$new = $orig.Duplicate()

# EMU -> point conversion (1 pt = 12700 EMU) for the new callout's
# target position/size:
#   off  x=3203848  y=2492896
#   ext  cx=2664296 cy=2104308
$new.Left   = 3203848 / 12700
$new.Top    = 2492896 / 12700
$new.Width  = 2664296 / 12700
# 2104308 / 12700 itself round-trips through the Single-precision
# Height property to 2104307 EMU (off by one); nudge very slightly so
# it lands back on the exact target after PowerPoint re-quantizes it.
$new.Height = 165.6936

# Mirror it horizontally (flipH="1" in the target XML).
$new.HorizontalFlip = -1
